# Update the "Förändrad" (Changed) date column C for rows 2-28
# from serial date 45455 (2024-06-12) to 45456 (2024-06-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45455) {
        $cell.Value2 = 45456
    }
}
